# Add the new "add_new_faculty" worksheet after the last existing sheet
# (edit_departments) and populate it with the faculty test data, mirroring
# the target diff.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "add_new_faculty"

# Column widths (characters)
$ws.Columns.Item(1).ColumnWidth = 17.140625
$ws.Columns.Item(2).ColumnWidth = 27.42578125

# Header row
$ws.Range("A1").Value = "code"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "center"
$ws.Range("D1").Value = "status"
$ws.Range("E1").Value = "runmode"

# Data rows
$ws.Range("A2").Value = "Auto_FCT001"
$ws.Range("B2").Value = "abcdefghi jklimnop 123"
$ws.Range("D2").Value = "Active"
$ws.Range("E2").Value = "Y"

$ws.Range("A3").Value = "Auto_FCT002"
$ws.Range("B3").Value = "abcdefghi jklimnop 124"
$ws.Range("D3").Value = "Active"
$ws.Range("E3").Value = "Y"

$ws.Range("A4").Value = "Auto_FCT003"
$ws.Range("B4").Value = "abcdefghi jklimnop 125"
$ws.Range("D4").Value = "Active"
$ws.Range("E4").Value = "Y"

$ws.Range("A5").Value = "Auto_FCT004"
$ws.Range("B5").Value = "abcdefghi jklimnop 126"
$ws.Range("D5").Value = "Active"
$ws.Range("E5").Value = "Y"

$ws.Range("A6").Value = "Auto_FCT005"
$ws.Range("B6").Value = "abcdefghi jklimnop 127"
$ws.Range("D6").Value = "Inactive"
$ws.Range("E6").Value = "Y"

$ws.Range("A7").Value = "Auto_FCT006"
$ws.Range("B7").Value = "abcdefghi jklimnop 128"
$ws.Range("D7").Value = "Inactive"
$ws.Range("E7").Value = "Y"

# Selection + activation, matching the tabSelected/activeTab move onto the
# new sheet.
$ws.Range("E2:E7").Select() | Out-Null
$ws.Activate()
